$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.304365634918213
$ws.Range("B1").Value = 2.085749626159668
$ws.Range("C1").Value = 4.71490478515625
$ws.Range("D1").Value = 3.469613313674927
$ws.Range("E1").Value = 1.356612682342529
